$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: the current week's rows (123:125) become last week's history.
# Copy them down to new rows 126:128 before overwriting them, preserving
# all formatting (e.g. the date style on column D).
$src = $ws.Range("A123:T125")
$dst = $ws.Range("A126:T128")
$src.Copy($dst)

# Step 2: update rows 123:125 in place with this week's figures.

# Row 123 - Especial
$ws.Cells.Item(123, 4).Value = 44516   # Fecha
$ws.Cells.Item(123, 14).Value = 26000  # Precio minimo
$ws.Cells.Item(123, 15).Value = 26000  # Precio maximo
$ws.Cells.Item(123, 16).Value = 26000  # Precio promedio ponderado
$ws.Cells.Item(123, 19).Value = 2600   # Precio $/Kg

# Row 124 - Primera
$ws.Cells.Item(124, 4).Value = 44516
$ws.Cells.Item(124, 13).Value = 50     # Volumen
$ws.Cells.Item(124, 14).Value = 24000
$ws.Cells.Item(124, 15).Value = 24000
$ws.Cells.Item(124, 16).Value = 24000
$ws.Cells.Item(124, 19).Value = 2400

# Row 125 - Segunda
$ws.Cells.Item(125, 4).Value = 44516
$ws.Cells.Item(125, 13).Value = 50     # Volumen (only field besides date that changed)

$wb.Save()
